$wb = $excel.ActiveWorkbook

# --- rs_contact_subject: fill in missing subject-contact pairings ---
$wsRcs = $wb.Worksheets.Item("rs_contact_subject")

# Row 2 currently has A2 = 'RS0' but B2 empty -> set subject S0
$wsRcs.Range("B2").Value = "S0"

# New rows 3-7
$wsRcs.Range("A3").Value = "RS0"
$wsRcs.Range("B3").Value = "S4"

$wsRcs.Range("A4").Value = "RS1"
$wsRcs.Range("B4").Value = "S4"

$wsRcs.Range("A5").Value = "RS1"
$wsRcs.Range("B5").Value = "S15"

$wsRcs.Range("A6").Value = "RS0"
$wsRcs.Range("B6").Value = "S25"

$wsRcs.Range("A7").Value = "RS0"
$wsRcs.Range("B7").Value = "S2"

$wsRcs.Activate() | Out-Null
$wsRcs.Range("A8").Select() | Out-Null

# --- ts_radiosilence: just a selection move ---
$wsTs = $wb.Worksheets.Item("ts_radiosilence")
$wsTs.Activate() | Out-Null
$wsTs.Range("E14").Select() | Out-Null

# --- rs_contact_location: selection move ---
$wsRcl = $wb.Worksheets.Item("rs_contact_location")
$wsRcl.Activate() | Out-Null
$wsRcl.Range("C3").Select() | Out-Null

# --- expenses: selection move (becomes non-active tab) ---
$wsExp = $wb.Worksheets.Item("expenses")
$wsExp.Activate() | Out-Null
$wsExp.Range("J2:J4").Select() | Out-Null

# --- locationList: selection move ---
$wsLoc = $wb.Worksheets.Item("locationList")
$wsLoc.Activate() | Out-Null
$wsLoc.Range("F2:F3").Select() | Out-Null

# --- subjectList: rename table reference in formula + selection, becomes active tab ---
$wsSubj = $wb.Worksheets.Item("subjectList")
$wsSubj.Range("C1").Formula = '="INSERT INTO subjectList (" & A1 & ", " &  B1 & ") values ("'
$wsSubj.Activate() | Out-Null
$wsSubj.Range("E24").Select() | Out-Null
